# Update the spelling of Turkey's country name in the OECD raw-data sheet,
# changing it to the current official spelling "Türkiye".
# (The previous spelling "Turkey" remains available as an alternate name
# in the shared-strings table once Excel resaves the workbook.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B36").Value = "Türkiye"
